$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "nguyen thi dong " record (originally row 5, id 8) -
# deleting shifts every row below it up by one, so the former row 6
# ("nguyen thi a ", id 9) becomes the new row 5.
$ws.Rows.Item(5).Delete()

# Remove the former row 7 ("nguyen thi b ", id 10), which is now row 6
# after the previous shift.
$ws.Rows.Item(6).Delete()
